$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "BACK END" header in A1 becomes a merged A1:D1 title cell reading
# "BREAKDOWN FOR BACK END", bold/size 18, centered, wrapped.
$ws.Range("A1").Value = "BREAKDOWN FOR BACK END"

$title = $ws.Range("A1:D1")
$title.Merge() | Out-Null
$title.Font.Bold = $true
$title.Font.Size = 18
$title.HorizontalAlignment = -4108
$title.VerticalAlignment = -4108
$title.WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Selection moves to A2 after the edit.
$ws.Range("A2").Select() | Out-Null
